$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.353.31"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "1.871.33"
$ws.Range("E3").Value = "  +0.33%  "
$ws.Range("D4").Value = "0.9976"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "0.7136"
$ws.Range("E5").Value = "  +1.44%  "
$ws.Range("D6").Value = "238.41"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").Value = "0.07913"
$ws.Range("E8").Value = "  -3.74%  "
$ws.Range("D9").Value = "0.3072"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("D10").Value = "25.30"
$ws.Range("E10").Value = "  +7.95%  "
$ws.Range("D11").Value = "0.08177"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "1.865.09"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").Value = "5.244"
$ws.Range("E13").Value = "  +0.89%  "
$ws.Range("D14").Value = "0.7221"
$ws.Range("E14").Value = "  +0.13%  "
$ws.Range("D15").Value = "89.41"
$ws.Range("E15").Value = "  -0.19%  "
$ws.Range("D16").Value = "29.406.94"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("D17").Value = "5.827"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "242.55"
$ws.Range("E18").Value = "  +1.80%  "
$ws.Range("D19").Value = "0.000007815"
$ws.Range("E19").Value = "  -1.15%  "
$ws.Range("D20").Value = "13.26"
$ws.Range("E20").Value = "  -1.38%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "2.107.87"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "0.9974"
$ws.Range("E23").Value = "  -0.34%  "
$ws.Range("D24").Value = "7.597"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").Value = "162.43"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "8.954"
$ws.Range("E26").Value = "  -0.71%  "
$ws.Range("D27").Value = "0.1458"
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").Value = "18.17"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").Value = "1.918"
$ws.Range("E29").Value = "  -3.93%  "
$ws.Range("D30").Value = "1.371"
$ws.Range("E30").Value = "  -4.62%  "
$ws.Range("D31").Value = "1.475"
$ws.Range("E31").Value = "  -0.81%  "
$ws.Range("D32").Value = "4.333"
$ws.Range("E32").Value = "  -2.58%  "
$ws.Range("D33").Value = "4.056"
$ws.Range("E33").Value = "  -0.43%  "
$ws.Range("D34").Value = "0.05210"
$ws.Range("E34").Value = "  -0.35%  "
$ws.Range("D35").Value = "1.188"
$ws.Range("E35").Value = "  +1.01%  "
$ws.Range("D36").Value = "0.7197"
$ws.Range("E36").Value = "  +1.74%  "
$ws.Range("D37").Value = "1.006"
$ws.Range("E37").Value = "  +0.03%  "
$ws.Range("E38").Value = "  +0.41%  "
$ws.Range("D39").Value = "0.01852"
$ws.Range("E39").Value = "  -0.24%  "
$ws.Range("D40").Value = "2.701"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D41").Value = "1.176.16"
$ws.Range("E41").Value = "  +2.20%  "
$ws.Range("D42").Value = "0.9190"
$ws.Range("E42").Value = "  -0.50%  "
$ws.Range("D43").Value = "6.006"
$ws.Range("E43").Value = "  +0.47%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "0.4289"
$ws.Range("E44").Value = "  -0.21%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "71.30"
$ws.Range("E45").Value = "  +0.29%  "
$ws.Range("D46").Value = "0.9999"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("D47").Value = "102.36"
$ws.Range("E47").Value = "  -0.79%  "
$ws.Range("D48").Value = "0.5351"
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("D49").Value = "1.757"
$ws.Range("E49").Value = "  -1.40%  "
$ws.Range("D50").Value = "9.171"
$ws.Range("E50").Value = "  -0.51%  "
$ws.Range("D51").Value = "7.013"
$ws.Range("E51").Value = "  +0.11%  "
